# Rename the "AddressBook" component to "TravelBuddy" on the Logic
# Component Class Diagram slide (matches commit: modify "addressbook",
# "address book", "person"/"persons" -> "travelbuddy", "place"/"places").
#
# The only shape on the slide containing the literal text "AddressBook"
# is the "Rectangle 62" box (shape Id 16) that together with the
# paragraph below it reads "AddressBookParser". Only the first
# paragraph ("AddressBook") needs to change; the second paragraph
# ("Parser") stays as-is.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$found = $false
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tf = $shp.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            if ($tr.Text.Contains("AddressBook")) {
                $tr.Paragraphs(1, 1).Text = "TravelBuddy"
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not locate the 'AddressBookParser' shape to rename"
}
